$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: reuse the common "pre-condition" boilerplate text (same value used in rows 10-21)
$ws.Range("C8").Value = "1. Account should be created in Tudunet`n2. Finaltest.ihex should be uploaded"

# Row 9: reuse the common boilerplate text for "pre-condition", "test case steps" and
# "expected result" columns (same values used in rows 10-21), and grow the row to fit.
$ws.Range("C9").Value = "1. Account should be created in Tudunet`n2. Finaltest.ihex should be uploaded"
$ws.Range("D9").Value = "1. Open web page ""http://www.tudunet.tu-darmstadt.de/""`n2. Enter ""username"" and ""password"".`n3. Select ""Jobs""`n4. Select ""Manage Jobs""`n5. Select ""Create new job""`n6. Write name ""TestFianl"" and description. `n7. Click on ""Next"" button and Select ""Finaltest.ihex"" file `n8. Click on ""Next"" button and  Select which program will run on individual nodes.`n9. Select ""playground"" nodes and click ""Finish"" button"
$ws.Range("E9").Value = "1. ""TestFinal"" job should be created."

$ws.Rows.Item(9).RowHeight = 240

# Clear the lingering cell selection saved in the sheet view
$ws.Range("A1").Select()
